# The workbook had a "laneNumber" column (column E) that is no longer
# part of the standard fastq-tracking template. Remove the entire column,
# which shifts every later column (sequencerModel, flowcellType, purpose,
# tapestationConc, volumePooled, readsObtained, fastqFileName, ...) one
# position to the left, and leaves the selection where the deleted column
# used to be.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("E").Delete()

# After a column delete, Excel leaves the active selection sitting on the
# column that used to be there (the cells that slid into its place).
$ws.Range("E1:E1048576").Select()
